# "ESPATULAS PLASTICAS DISMAY" price list - fix gui step 1 and 2:
#   1) bump the printed date in A1 by one day (45308 -> 45309)
#   2) double the unit prices for the two items (D29, D30)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: update the date shown at the top of the sheet.
$ws.Range("A1").Value = 45309

# Step 2: update the two unit prices (both simply doubled).
$ws.Range("D29").Value = 112.4
$ws.Range("D30").Value = 187.2
